$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = -7.906999999999999
$ws.Range("D12").Value = -7.885999999999998
$ws.Range("E13").Value = 12.583
$ws.Range("D18").Value = -7.952000000000001
